$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.299.05'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '2.425.02'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.08%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.39'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +3.64%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.39'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +1.90%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("E8").Value = '  -0.09%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +4.92%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.64'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E12").Value = '  -1.12%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.21'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -2.59%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.97'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '2.804.87'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = '2.440.03'
$ws.Range("E16").Value = '  -0.18%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.832'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '45.234.09'
$ws.Range("E18").Value = '  +1.69%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.27'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -1.23%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.36'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").Value = '0.0₃0926'
$ws.Range("E21").Value = '  +2.34%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.04'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +3.28%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.39'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +1.82%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.72'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +2.07%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  +0.04%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.66'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +1.73%  '
$ws.Range("E28").Value = '  -0.74%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.63'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  -0.15%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.42'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +0.82%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.10'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.127'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +5.85%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.33'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +4.89%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.24'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("E35").Value = '  +0.19%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0755'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("E37").Value = '  -0.22%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '128.58'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.90'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  +1.21%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.75'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -4.59%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0290'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").Value = '1.948.74'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("E46").Value = '  -2.78%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.95'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.79'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +7.52%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -3.76%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.44'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +3.75%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.80'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +3.89%  '
